$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new header columns: Email (D) and Website (E)
$ws.Range("D1").Value = "Email"
$ws.Range("E1").Value = "Website"

# Add a mailto hyperlink for Don's email in D3
$ws.Hyperlinks.Add($ws.Range("D3"), "mailto:greg.smalter@gmail.com", "", "", "greg.smalter@gmail.com") | Out-Null

# Add a web hyperlink for Eddie's website in E6
$ws.Hyperlinks.Add($ws.Range("E6"), "https://www.google.com", "", "", "https://www.google.com") | Out-Null

# Leave the selection where the author ended up editing
$ws.Range("E7").Select() | Out-Null
